$d = $word.ActiveDocument

$replacements = @(
    @{old = "62÷5="; new = "17÷8="},
    @{old = "50÷4="; new = "13÷9="},
    @{old = "97÷4="; new = "57÷8="},
    @{old = "55÷3="; new = "69÷2="},
    @{old = "11÷7="; new = "11÷4="},
    @{old = "62÷7="; new = "31÷3="},
    @{old = "53÷5="; new = "67÷8="},
    @{old = "76÷4="; new = "89÷9="},
    @{old = "55÷9="; new = "68÷2="},
    @{old = "12÷7="; new = "17÷2="},
    @{old = "48÷7="; new = "98÷2="},
    @{old = "56÷2="; new = "20÷6="},
    @{old = "62÷8="; new = "64÷5="},
    @{old = "22÷7="; new = "50÷3="},
    @{old = "73÷7="; new = "21÷6="},
    @{old = "19÷2="; new = "45÷8="},
    @{old = "87÷5="; new = "85÷2="},
    @{old = "41÷7="; new = "10÷2="},
    @{old = "42÷3="; new = "70÷3="},
    @{old = "11÷3="; new = "55÷4="},
    @{old = "12÷8="; new = "88÷9="},
    @{old = "76÷6="; new = "97÷7="},
    @{old = "68÷3="; new = "57÷7="},
    @{old = "75÷9="; new = "22÷5="},
    @{old = "62÷9="; new = "14÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
